$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'69.611.43"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -0.77%  "
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'3.838.04"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  +2.65%  "
$ws.Cells.Item(3, 5).Style = "Normal"

$ws.Cells.Item(4, 5).Value = "'  +0.03%  "
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'607.32"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  -2.56%  "
$ws.Cells.Item(5, 5).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'173.64"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -3.85%  "
$ws.Cells.Item(6, 5).Style = "Normal"

$ws.Cells.Item(7, 4).Value = "'3.836.89"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  +2.66%  "
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 5).Value = "'  +0.08%  "
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 4).Value = "'0.525"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -1.91%  "
$ws.Cells.Item(9, 5).Style = "Normal"

$ws.Cells.Item(10, 4).Value = "'0.165"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -1.46%  "
$ws.Cells.Item(10, 5).Style = "Normal"

$ws.Cells.Item(11, 5).Value = "'  +2.57%  "
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 5).Value = "'  -1.73%  "
$ws.Cells.Item(12, 5).Style = "Normal"

$ws.Cells.Item(13, 4).Value = "'39.74"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -2.83%  "
$ws.Cells.Item(13, 5).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'0.0000253"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -2.50%  "
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 4).Value = "'4.480.03"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +2.69%  "
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "'3.836.87"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +2.69%  "
$ws.Cells.Item(16, 5).Style = "Normal"

$ws.Cells.Item(17, 4).Value = "'69.765.05"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  -0.58%  "
$ws.Cells.Item(17, 5).Style = "Normal"

$ws.Cells.Item(18, 4).Value = "'7.43"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -2.30%  "
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 5).Value = "'  -3.24%  "
$ws.Cells.Item(19, 5).Style = "Normal"

$ws.Cells.Item(20, 4).Value = "'16.55"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -1.89%  "
$ws.Cells.Item(20, 5).Style = "Normal"

$ws.Cells.Item(21, 4).Value = "'503.68"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  -0.41%  "
$ws.Cells.Item(21, 5).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'9.45"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  +1.13%  "
$ws.Cells.Item(22, 5).Style = "Normal"

$ws.Cells.Item(23, 5).Value = "'  +1.46%  "
$ws.Cells.Item(23, 5).Style = "Normal"

$ws.Cells.Item(24, 2).Value = "'Litecoin"
$ws.Cells.Item(24, 2).Style = "Normal"
$ws.Cells.Item(24, 3).Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(24, 3).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'85.77"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -1.08%  "
$ws.Cells.Item(24, 5).Style = "Normal"

$ws.Cells.Item(25, 2).Value = "'Fetch.AI"
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 3).Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(25, 3).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'2.44"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -5.59%  "
$ws.Cells.Item(25, 5).Style = "Normal"

$ws.Cells.Item(26, 5).Value = "'  +4.95%  "
$ws.Cells.Item(26, 5).Style = "Normal"

$ws.Cells.Item(27, 4).Value = "'12.57"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  -4.69%  "
$ws.Cells.Item(27, 5).Style = "Normal"

$ws.Cells.Item(28, 4).Value = "'10.35"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  -10.00%  "
$ws.Cells.Item(28, 5).Style = "Normal"

$ws.Cells.Item(29, 5).Value = "'  +0.19%  "
$ws.Cells.Item(29, 5).Style = "Normal"

$ws.Cells.Item(30, 4).Value = "'2.51"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  +0.52%  "
$ws.Cells.Item(30, 5).Style = "Normal"

$ws.Cells.Item(31, 4).Value = "'2.97"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  +0.66%  "
$ws.Cells.Item(31, 5).Style = "Normal"

$ws.Cells.Item(32, 4).Value = "'7.92"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -0.74%  "
$ws.Cells.Item(32, 5).Style = "Normal"

$ws.Cells.Item(33, 4).Value = "'32.15"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +2.95%  "
$ws.Cells.Item(33, 5).Style = "Normal"

$ws.Cells.Item(34, 4).Value = "'0.112"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -2.65%  "
$ws.Cells.Item(34, 5).Style = "Normal"

$ws.Cells.Item(35, 5).Value = "'  -0.12%  "
$ws.Cells.Item(35, 5).Style = "Normal"

$ws.Cells.Item(36, 5).Value = "'  -1.93%  "
$ws.Cells.Item(36, 5).Style = "Normal"

$ws.Cells.Item(37, 4).Value = "'6.05"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -2.09%  "
$ws.Cells.Item(37, 5).Style = "Normal"

$ws.Cells.Item(38, 4).Value = "'0.140"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +2.39%  "
$ws.Cells.Item(38, 5).Style = "Normal"

$ws.Cells.Item(39, 4).Value = "'483.15"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +12.25%  "
$ws.Cells.Item(39, 5).Style = "Normal"

$ws.Cells.Item(40, 4).Value = "'0.334"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -1.24%  "
$ws.Cells.Item(40, 5).Style = "Normal"

$ws.Cells.Item(41, 2).Value = "'OKB"
$ws.Cells.Item(41, 2).Style = "Normal"
$ws.Cells.Item(41, 3).Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(41, 3).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'49.67"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  -1.41%  "
$ws.Cells.Item(41, 5).Style = "Normal"

$ws.Cells.Item(42, 2).Value = "'Stacks"
$ws.Cells.Item(42, 2).Style = "Normal"
$ws.Cells.Item(42, 3).Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 3).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'2.04"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  -3.12%  "
$ws.Cells.Item(42, 5).Style = "Normal"

$ws.Cells.Item(43, 5).Value = "'  +3.65%  "
$ws.Cells.Item(43, 5).Style = "Normal"

$ws.Cells.Item(44, 4).Value = "'43.06"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -6.21%  "
$ws.Cells.Item(44, 5).Style = "Normal"

$ws.Cells.Item(45, 4).Value = "'8.49"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -2.89%  "
$ws.Cells.Item(45, 5).Style = "Normal"

$ws.Cells.Item(46, 4).Value = "'2.922.23"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -2.94%  "
$ws.Cells.Item(46, 5).Style = "Normal"

$ws.Cells.Item(47, 4).Value = "'0.0359"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -1.55%  "
$ws.Cells.Item(47, 5).Style = "Normal"

$ws.Cells.Item(48, 4).Value = "'139.90"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +1.63%  "
$ws.Cells.Item(48, 5).Style = "Normal"

$ws.Cells.Item(49, 5).Value = "'  +0.06%  "
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(50, 4).Value = "'26.63"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  -3.71%  "
$ws.Cells.Item(50, 5).Style = "Normal"

$ws.Cells.Item(51, 5).Value = "'  -4.54%  "
$ws.Cells.Item(51, 5).Style = "Normal"
